$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E3").Value = 16.2242
$ws.Range("D12").Value = -7.102899999999996
$ws.Range("E14").Value = 16.67730000000001
$ws.Range("E26").Value = 16.29409999999999
$ws.Range("E31").Value = 16.6588
$ws.Range("D32").Value = -8.797699999999999
$ws.Range("E35").Value = 16.6609
$ws.Range("D36").Value = -8.228899999999996
$ws.Range("E37").Value = 16.5807
$ws.Range("D38").Value = -7.783499999999998
$ws.Range("E45").Value = 16.5928
$ws.Range("D46").Value = -8.3149
$ws.Range("D54").Value = -8.174300000000001
$ws.Range("D55").Value = -8.286799999999998
$ws.Range("E57").Value = 16.6523
$ws.Range("D67").Value = -6.2532
$ws.Range("D69").Value = -7.052099999999996
$ws.Range("D72").Value = -7.4337
$ws.Range("D91").Value = -6.641699999999998
$ws.Range("D99").Value = -7.751199999999996
$ws.Range("E100").Value = 16.41600000000001
$ws.Range("E102").Value = 16.72469999999999
